{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2025-07-17 Thursday\", \"2025-07-18 Friday\"],\n  [\"42\u00d754=2268\", \"19\u00d747=893\"],\n  [\"11\u00d778=858\", \"57\u00d752=2964\"],\n  [\"16\u00d713=208\", \"25\u00d751=1275\"],\n  [\"22\u00d788=1936\", \"55\u00d795=5225\"],\n  [\"41\u00d760=2460\", \"24\u00d766=1584\"],\n  [\"78\u00d714=1092\", \"28\u00d787=2436\"],\n  [\"90\u00d724=2160\", \"96\u00d798=9408\"],\n  [\"57\u00d726=1482\", \"59\u00d776=4484\"],\n  [\"37\u00d769=2553\", \"92\u00d757=5244\"],\n  [\"96\u00d769=6624\", \"69\u00d758=4002\"],\n  [\"54\u00d741=2214\", \"95\u00d771=6745\"],\n  [\"32\u00d737=1184\", \"24\u00d761=1464\"],\n  [\"97\u00d797=9409\", \"35\u00d779=2765\"],\n  [\"92\u00d727=2484\", \"69\u00d718=1242\"],\n  [\"90\u00d796=8640\", \"18\u00d767=1206\"],\n  [\"11\u00d771=781\", \"23\u00d716=368\"],\n  [\"87\u00d731=2697\", \"93\u00d742=3906\"],\n  [\"52\u00d770=3640\", \"38\u00d730=1140\"],\n  [\"70\u00d711=770\", \"35\u00d785=2975\"],\n  [\"53\u00d778=4134\", \"36\u00d791=3276\"],\n  [\"85\u00d737=3145\", \"63\u00d717=1071\"],\n  [\"40\u00d759=2360\", \"34\u00d764=2176\"],\n  [\"38\u00d725=950\", \"83\u00d761=5063\"],\n  [\"19\u00d746=874\", \"92\u00d772=6624\"],\n  [\"91\u00d774=6734\", \"81\u00d726=2106\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-17 Thursday\", \"2025-07-18 Friday\"),\n    @(\"42\u00d754=2268\", \"19\u00d747=893\"),\n    @(\"11\u00d778=858\", \"57\u00d752=2964\"),\n    @(\"16\u00d713=208\", \"25\u00d751=1275\"),\n    @(\"22\u00d788=1936\", \"55\u00d795=5225\"),\n    @(\"41\u00d760=2460\", \"24\u00d766=1584\"),\n    @(\"78\u00d714=1092\", \"28\u00d787=2436\"),\n    @(\"90\u00d724=2160\", \"96\u00d798=9408\"),\n    @(\"57\u00d726=1482\", \"59\u00d776=4484\"),\n    @(\"37\u00d769=2553\", \"92\u00d757=5244\"),\n    @(\"96\u00d769=6624\", \"69\u00d758=4002\"),\n    @(\"54\u00d741=2214\", \"95\u00d771=6745\"),\n    @(\"32\u00d737=1184\", \"24\u00d761=1464\"),\n    @(\"97\u00d797=9409\", \"35\u00d779=2765\"),\n    @(\"92\u00d727=2484\", \"69\u00d718=1242\"),\n    @(\"90\u00d796=8640\", \"18\u00d767=1206\"),\n    @(\"11\u00d771=781\", \"23\u00d716=368\"),\n    @(\"87\u00d731=2697\", \"93\u00d742=3906\"),\n    @(\"52\u00d770=3640\", \"38\u00d730=1140\"),\n    @(\"70\u00d711=770\", \"35\u00d785=2975\"),\n    @(\"53\u00d778=4134\", \"36\u00d791=3276\"),\n    @(\"85\u00d737=3145\", \"63\u00d717=1071\"),\n    @(\"40\u00d759=2360\", \"34\u00d764=2176\"),\n    @(\"38\u00d725=950\", \"83\u00d761=5063\"),\n    @(\"19\u00d746=874\", \"92\u00d772=6624\"),\n    @(\"91\u00d774=6734\", \"81\u00d726=2106\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
